# Updates the cryptos list (prices in column D, 1h-volume % in column E)
# to reflect the latest scrape, as produced by the GitHub Actions job on
# Mon Mar 27 07:07:37 UTC 2023.
#
# Notes on technique:
#  - Column D ("Price") values are stored as plain text in the workbook
#    (e.g. "27.834.54", "0.5802"), even though many of them look like
#    numbers. Assigning a numeric-looking string straight to .Value makes
#    Excel auto-convert it to a real number (dropping trailing zeros,
#    "0.5800" -> 0.58, etc.), which would not match the source data. To
#    keep these values as text we prefix them with a leading apostrophe,
#    exactly like a user typing a text-forced value into Excel.
#  - Column E ("Volume(1h)") values are already padded with spaces
#    (e.g. "  +1.11%  "), so Excel always treats them as plain text and a
#    normal assignment is sufficient.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'27.848.85"
$ws.Range("E2").Value = "  +1.11%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.757.25"
$ws.Range("E3").Value = "  +0.16%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'327.76"
$ws.Range("E5").Value = "  +1.05%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.02%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.4555"
$ws.Range("E7").Value = "  +0.14%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3494"
$ws.Range("E8").Value = "  -1.47%  "

# Row 9 - OKB
$ws.Range("D9").Value = "'41.94"
$ws.Range("E9").Value = "  +1.09%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.07355"
$ws.Range("E10").Value = "  -1.50%  "

# Row 11 - Polygon
$ws.Range("E11").Value = "  -0.15%  "

# Row 12 - BinanceUSD
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.02%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'5.977"
$ws.Range("E14").Value = "  -0.54%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "'7.174"
$ws.Range("E15").Value = "  +0.17%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "'1.756.71"
$ws.Range("E16").Value = "  -0.49%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "'91.61"
$ws.Range("E17").Value = "  -2.24%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "'0.00001054"
$ws.Range("E18").Value = "  +0.06%  "

# Row 19 - TRON
$ws.Range("D19").Value = "'0.06409"
$ws.Range("E19").Value = "  +0.31%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  +0.05%  "

# Row 21 - Avalanche
$ws.Range("D21").Value = "'16.84"
$ws.Range("E21").Value = "  -1.53%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'5.754"
$ws.Range("E22").Value = "  +0.33%  "

# Row 23 - WrappedBTC
$ws.Range("D23").Value = "'27.880.16"
$ws.Range("E23").Value = "  +1.03%  "

# Row 24 - Cosmos
$ws.Range("D24").Value = "'11.17"
$ws.Range("E24").Value = "  -0.26%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +4.06%  "

# Row 26 - Monero
$ws.Range("D26").Value = "'162.24"
$ws.Range("E26").Value = "  -1.92%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'20.04"
$ws.Range("E27").Value = "  -0.43%  "

# Row 28 - WrappedliquidstakedEther2.0
$ws.Range("D28").Value = "'1.960.28"
$ws.Range("E28").Value = "  -0.12%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "'2.162"
$ws.Range("E29").Value = "  +1.43%  "

# Row 30 - BitcoinCash
$ws.Range("D30").Value = "'123.44"
$ws.Range("E30").Value = "  -1.73%  "

# Row 31 - ImmutableX
$ws.Range("D31").Value = "'1.075"
$ws.Range("E31").Value = "  -1.24%  "

# Row 32 - Stellar
$ws.Range("D32").Value = "'0.09262"
$ws.Range("E32").Value = "  +0.59%  "

# Row 33 - HuobiToken
$ws.Range("D33").Value = "'3.660"
$ws.Range("E33").Value = "  -0.01%  "

# Row 34 - Filecoin
$ws.Range("D34").Value = "'5.532"
$ws.Range("E34").Value = "  +0.00%  "

# Row 35 - Aptos
$ws.Range("D35").Value = "'11.74"
$ws.Range("E35").Value = "  +0.23%  "

# Row 36 - Hedera
$ws.Range("D36").Value = "'0.06105"
$ws.Range("E36").Value = "  +1.50%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "'0.02261"
$ws.Range("E37").Value = "  -0.92%  "

# Row 38 - Algorand
$ws.Range("D38").Value = "'0.2062"
$ws.Range("E38").Value = "  -1.70%  "

# Row 39 - InternetComputer(DFINITY)
$ws.Range("D39").Value = "'4.895"
$ws.Range("E39").Value = "  -0.45%  "

# Row 40 - TheSandbox
$ws.Range("D40").Value = "'0.6186"
$ws.Range("E40").Value = "  -1.56%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "'1.178"
$ws.Range("E41").Value = "  -0.16%  "

# Row 42 - WEMIXTOKEN
$ws.Range("E42").Value = "  -1.46%  "

# Row 43 - FraxShare
$ws.Range("D43").Value = "'7.773"
$ws.Range("E43").Value = "  -0.50%  "

# Row 44 - EnergySwap
$ws.Range("D44").Value = "'13.10"
$ws.Range("E44").Value = "  -0.19%  "

# Row 45 - PancakeSwap
$ws.Range("D45").Value = "'3.727"
$ws.Range("E45").Value = "  +0.29%  "

# Row 46 - Decentraland
$ws.Range("D46").Value = "'0.5800"

# Row 47 - Quant
$ws.Range("D47").Value = "'122.53"
$ws.Range("E47").Value = "  +0.33%  "

# Row 48 - NEARProtocol
$ws.Range("D48").Value = "'1.928"
$ws.Range("E48").Value = "  -0.26%  "

# Row 49 - was Cronos, now EOS (rows 49/50 swapped order)
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "'1.121"
$ws.Range("E49").Value = "  -0.78%  "

# Row 50 - was EOS, now Cronos
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.06787"
$ws.Range("E50").Value = "  -1.54%  "

# Row 51 - Aave
$ws.Range("D51").Value = "'72.18"
$ws.Range("E51").Value = "  -0.01%  "
